# "Actualización automática desde Jupyter"
# The route_id column (A2:A130) changes from "JQ_R" to "80JQ_R" for every
# bus-stop row, and the sheet's view/selection moves from the bottom of the
# list (C119) up to the top of the route_id column (A2:A130, scrolled near A7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update every route_id cell in the data body (rows 2-130) to the new code.
$ws.Range("A2:A130").Value = "80JQ_R"

# Move the visible selection to the route_id column, matching where the
# author was working when the workbook was re-saved.
$ws.Activate()
$ws.Range("A2:A130").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
